# Apply the cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values scraped for this run (column letter + row -> new text)
$newValues = @{
    "D2" = '28.131.28'
    "E2" = '  -1.27%  '
    "D3" = '1.792.57'
    "E3" = '  -1.46%  '
    "E4" = '  +0.12%  '
    "D5" = '313.79'
    "E5" = '  -0.35%  '
    "E6" = '  +0.15%  '
    "E7" = '  +1.68%  '
    "D8" = '0.3807'
    "E8" = '  -3.40%  '
    "D9" = '0.08004'
    "E9" = '  -2.71%  '
    "D10" = '41.44'
    "E10" = '  -0.45%  '
    "D11" = '1.094'
    "E11" = '  -1.30%  '
    "D12" = '6.280'
    "E12" = '  -0.55%  '
    "D13" = '1.003'
    "E13" = '  +0.14%  '
    "D14" = '20.48'
    "E14" = '  -2.98%  '
    "D15" = '7.287'
    "E15" = '  -3.27%  '
    "D16" = '1.784.06'
    "E16" = '  -1.74%  '
    "D17" = '91.54'
    "E17" = '  -0.99%  '
    "D18" = '0.00001092'
    "E18" = '  -3.62%  '
    "D19" = '0.06564'
    "E19" = '  -1.31%  '
    "D20" = '1.002'
    "E20" = '  +0.19%  '
    "D21" = '17.30'
    "E21" = '  -2.68%  '
    "D22" = '5.949'
    "E22" = '  -2.28%  '
    "D23" = '28.169.50'
    "E23" = '  -1.27%  '
    "D24" = '11.13'
    "E24" = '  -2.60%  '
    "D25" = '2.273'
    "E25" = '  +0.28%  '
    "D26" = '160.41'
    "E26" = '  +2.98%  '
    "D27" = '20.43'
    "E27" = '  -3.85%  '
    "D28" = '1.993.65'
    "E28" = '  -1.59%  '
    "D29" = '2.327'
    "E29" = '  -3.07%  '
    "D30" = '122.87'
    "E30" = '  -2.17%  '
    "D31" = '0.1077'
    "E31" = '  -1.63%  '
    "D32" = '1.053'
    "E32" = '  -5.15%  '
    "E33" = '  +0.37%  '
    "D34" = '5.540'
    "E34" = '  -4.07%  '
    "D35" = '0.07180'
    "E35" = '  +1.42%  '
    "D36" = '12.00'
    "E36" = '  +6.64%  '
    "D37" = '0.02306'
    "E37" = '  -1.65%  '
    "D38" = '0.2145'
    "E38" = '  -3.39%  '
    "D39" = '5.056'
    "E39" = '  -3.27%  '
    "D40" = '8.631'
    "E40" = '  -2.15%  '
    "D41" = '0.6160'
    "E41" = '  -2.32%  '
    "D42" = '1.163'
    "E42" = '  -1.31%  '
    "D43" = '13.21'
    "E43" = '  -1.88%  '
    "D47" = '127.51'
    "E47" = '  +2.14%  '
    "E48" = '  +2.55%  '
    "D49" = '1.918'
    "E49" = '  -3.37%  '
    "D50" = '0.06753'
    "E50" = '  -2.00%  '
    "D51" = '72.58'
    "E51" = '  -2.32%  '
    "B44" = 'WEMIXTOKEN'
    "C44" = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    "D44" = '1.309'
    "E44" = '  -6.55%  '
    "B45" = 'PancakeSwap'
    "C45" = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    "D45" = '3.761'
    "E45" = '  +0.78%  '
    "B46" = 'Decentraland'
    "C46" = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    "D46" = '0.5958'
    "E46" = '  +0.71%  '
}

# Cells whose new value looks like a plain number (e.g. "313.79") must stay
# text, matching the inline-string cells already used throughout column D.
$textForceCells = @(
    "D5", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19",
    "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D34",
    "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D47", "D49", "D50",
    "D51", "D44", "D45", "D46"
)

foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

foreach ($cell in $newValues.Keys) {
    $ws.Range($cell).Value = $newValues[$cell]
}

